# New Features and Bug Fixes
# Applies updates to the BugList workbook (Sheet1):
#   - Marks a couple of existing bugs as "done"
#   - Adds two new bug rows (one marked DONE, one marked done)
#   - Updates the saved view/selection to where the user left off

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark existing "MODULE TIDAK AKTIF TIDAK DITAMPILKAN" bug (row 12) as done
$ws.Range("C12").Value = "done"

# New bug row: perhitungan total pada print out salah dgn transaksi kredit salah (row 65) - done
$ws.Range("B65").Value = "perhitungan total pada print out salah dgn transaksi kredit salah"
$ws.Range("C65").Value = "done"

# New bug row: kolom status seharusnya tidak bisa diedit (row 56) - already DONE
$ws.Range("B56").Value = "kolom status seharusnya tidak bisa diedit"
$ws.Range("C56").Value = "DONE"

# Mark existing customer-group bug (row 58) as done
$ws.Range("C58").Value = "done"

# Update the view to reflect where the user scrolled/selected last
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("C57").Select()
